# Generate Report for Handback
# Adds a new handback row (e3be380e-ae4f-4390-af76-5dc4d5c584e1.md) to the
# Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$fileName   = "e3be380e-ae4f-4390-af76-5dc4d5c584e1.md"
$pathName   = "e2e\e3be380e-ae4f-4390-af76-5dc4d5c584e1.md"
$ext        = ".md"
$statusSync = "Handed back: in sync with en-US"
$genDate    = "2016-08-31 16:51:41"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()
$rOverview = $rowOverview.Range.Row

$wsOverview.Cells.Item($rOverview, 1).Value = $fileName
$wsOverview.Cells.Item($rOverview, 3).Value = $ext
$wsOverview.Cells.Item($rOverview, 5).Value = $statusSync
$wsOverview.Cells.Item($rOverview, 6).Value = $statusSync
$wsOverview.Cells.Item($rOverview, 7).Value = $genDate

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($rOverview, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3be380e0ae4f4390af765dc4d5c584e1e3be38/e2e/e3be380e-ae4f-4390-af76-5dc4d5c584e1.md",
    "",
    "",
    $pathName
) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$rowZhCn = $loZhCn.ListRows.Add()
$rZhCn = $rowZhCn.Range.Row

$zhXliff = "e3be380e-ae4f-4390-af76-5dc4d5c584e1.92cfe971ee3fefd0081b51016e2afa2520d651aa.zh-cn.xlf"

$wsZhCn.Cells.Item($rZhCn, 2).Value  = $ext
$wsZhCn.Cells.Item($rZhCn, 3).Value  = $statusSync
$wsZhCn.Cells.Item($rZhCn, 4).Value  = "e2e"
$wsZhCn.Cells.Item($rZhCn, 5).Value  = "ht"
$wsZhCn.Cells.Item($rZhCn, 6).Value  = "True"
$wsZhCn.Cells.Item($rZhCn, 7).Value  = $zhXliff
$wsZhCn.Cells.Item($rZhCn, 8).Value  = "2016-08-31 16:51:37"
$wsZhCn.Cells.Item($rZhCn, 10).Value = $zhXliff
$wsZhCn.Cells.Item($rZhCn, 11).Value = "2016-08-31 16:51:56"
$wsZhCn.Cells.Item($rZhCn, 13).Value = "True"
$wsZhCn.Cells.Item($rZhCn, 15).Value = "False"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item($rZhCn, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3be380e0ae4f4390af765dc4d5c584e1e3be38/e2e/e3be380e-ae4f-4390-af76-5dc4d5c584e1.md",
    "",
    "",
    $fileName
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item($rZhCn, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/92cfe971ee3fefd0081b51016e2afa2520d651a/e2e/e3be380e-ae4f-4390-af76-5dc4d5c584e1.md",
    "",
    "",
    $fileName
) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$rowDeDe = $loDeDe.ListRows.Add()
$rDeDe = $rowDeDe.Range.Row

$deXliff = "e3be380e-ae4f-4390-af76-5dc4d5c584e1.92cfe971ee3fefd0081b51016e2afa2520d651aa.de-de.xlf"

$wsDeDe.Cells.Item($rDeDe, 2).Value  = $ext
$wsDeDe.Cells.Item($rDeDe, 3).Value  = $statusSync
$wsDeDe.Cells.Item($rDeDe, 4).Value  = "e2e"
$wsDeDe.Cells.Item($rDeDe, 5).Value  = "ht"
$wsDeDe.Cells.Item($rDeDe, 6).Value  = "True"
$wsDeDe.Cells.Item($rDeDe, 7).Value  = $deXliff
$wsDeDe.Cells.Item($rDeDe, 8).Value  = $genDate
$wsDeDe.Cells.Item($rDeDe, 10).Value = $deXliff
$wsDeDe.Cells.Item($rDeDe, 11).Value = "2016-08-31 16:52:13"
$wsDeDe.Cells.Item($rDeDe, 13).Value = "True"
$wsDeDe.Cells.Item($rDeDe, 15).Value = "False"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item($rDeDe, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3be380e0ae4f4390af765dc4d5c584e1e3be38/e2e/e3be380e-ae4f-4390-af76-5dc4d5c584e1.md",
    "",
    "",
    $fileName
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item($rDeDe, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/92cfe971ee3fefd0081b51016e2afa2520d651a/e2e/e3be380e-ae4f-4390-af76-5dc4d5c584e1.md",
    "",
    "",
    $fileName
) | Out-Null
